# Applies the odds updates described in the commit "Atualizando o arquivo XLSX"
# (updated Betfair back/lay odds values for several matches on Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("S2").Value = 3.75
$ws.Range("AI2").Value = 40
$ws.Range("AL2").Value = 1000

# Row 3
$ws.Range("N3").Value = 1.37
$ws.Range("O3").Value = 1.23
$ws.Range("Q3").Value = 1.62
$ws.Range("S3").Value = 1.63

# Row 4
$ws.Range("F4").Value = 4.3
$ws.Range("G4").Value = 6
$ws.Range("H4").Value = 1.74
$ws.Range("I4").Value = 1.9
$ws.Range("J4").Value = 3.4
$ws.Range("K4").Value = 4.6
$ws.Range("L4").Value = 1.27
$ws.Range("M4").Value = 1.05
$ws.Range("N4").Value = 4.5
$ws.Range("O4").Value = 1.23
$ws.Range("P4").Value = 2.18
$ws.Range("Q4").Value = 1.6
$ws.Range("R4").Value = 1.49
$ws.Range("S4").Value = 2.6
$ws.Range("T4").Value = 1.7
$ws.Range("U4").Value = 2.2
$ws.Range("V4").Value = 2.1
$ws.Range("W4").Value = 1.22
$ws.Range("Y4").Value = 1000
$ws.Range("AC4").Value = 990
$ws.Range("AO4").Value = 11

# Row 5
$ws.Range("F5").Value = 3.8
$ws.Range("G5").Value = 4.7
$ws.Range("H5").Value = 1.93
$ws.Range("I5").Value = 2.14
$ws.Range("J5").Value = 3.35
$ws.Range("N5").Value = 3.5
$ws.Range("P5").Value = 1.86
$ws.Range("Q5").Value = 1.9
$ws.Range("S5").Value = 3.3
$ws.Range("V5").Value = 1.87
$ws.Range("AC5").Value = 9.8

# Row 7
$ws.Range("F7").Value = 7.4
$ws.Range("G7").Value = 9.6
$ws.Range("I7").Value = 1.5
$ws.Range("J7").Value = 4.6
$ws.Range("K7").Value = 5.9
$ws.Range("N7").Value = 4.3
$ws.Range("P7").Value = 2.1
$ws.Range("Q7").Value = 1.72
$ws.Range("R7").Value = 1.44
$ws.Range("U7").Value = 1.84
$ws.Range("V7").Value = 2.96
$ws.Range("W7").Value = 1.12

# Row 8
$ws.Range("N8").Value = 5.8

# Row 9
$ws.Range("I9").Value = 3.7
$ws.Range("J9").Value = 3.8
$ws.Range("V9").Value = 1.38
$ws.Range("X9").Value = 26
$ws.Range("Y9").Value = 22
$ws.Range("AA9").Value = 60
$ws.Range("AB9").Value = 17
$ws.Range("AF9").Value = 20
$ws.Range("AG9").Value = 13.5
$ws.Range("AH9").Value = 18
$ws.Range("AK9").Value = 22
$ws.Range("AN9").Value = 13

# Row 10
$ws.Range("I10").Value = 1.74
$ws.Range("V10").Value = 2.34

# Row 12
$ws.Range("G12").Value = 5.5
$ws.Range("H12").Value = 1.67
$ws.Range("O12").Value = 1.21
$ws.Range("S12").Value = 2.52
$ws.Range("T12").Value = 1.66
$ws.Range("U12").Value = 2.24

# Row 13
$ws.Range("K13").Value = 7.4
$ws.Range("P13").Value = 2.62
$ws.Range("R13").Value = 1.65
$ws.Range("U13").Value = 1.85

Write-Output "Applied all changes"